# Daily refresh of the cryptos.xlsx price/volume snapshot (GitHub Actions job).
# "Updated symbol list on Tue Jan  3 10:29:22 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every data cell on this sheet is stored as literal text even though Price/
# Volume(1h) look numeric (e.g. "246.16", "-0.31%"). Prefixing the assignment
# with an apostrophe forces Excel to keep it as text instead of re-typing it
# as a Number/Percentage, matching how the sheet was originally authored.
$q = "'"

$ws.Range("D2").Value = $q + '245.62'
$ws.Range("E2").Value = $q + '-0.46%'
$ws.Range("D3").Value = $q + '30.27'
$ws.Range("E3").Value = $q + '0.56%'
$ws.Range("D4").Value = $q + '5.150'
$ws.Range("E4").Value = $q + '-0.53%'
$ws.Range("D6").Value = $q + '6.657'
$ws.Range("E6").Value = $q + '0.83%'
$ws.Range("D7").Value = $q + '3.248'
$ws.Range("E7").Value = $q + '6.12%'
$ws.Range("D8").Value = $q + '0.8493'
$ws.Range("E8").Value = $q + '-1.41%'
$ws.Range("D9").Value = $q + '0.8549'
$ws.Range("E9").Value = $q + '-2.77%'
$ws.Range("D10").Value = $q + '0.1392'
$ws.Range("E10").Value = $q + '2.11%'
$ws.Range("D11").Value = $q + '0.07077'
$ws.Range("E11").Value = $q + '-0.07%'
$ws.Range("D12").Value = $q + '0.03257'
$ws.Range("E12").Value = $q + '12.73%'
$ws.Range("D13").Value = $q + '0.09373'
$ws.Range("E13").Value = $q + '-0.33%'
$ws.Range("D14").Value = $q + '0.001525'
$ws.Range("E14").Value = $q + '0.30%'
$ws.Range("D15").Value = $q + '0.0005958'
$ws.Range("D16").Value = $q + '0.005989'
$ws.Range("E16").Value = $q + '0.42%'
$ws.Range("D17").Value = $q + '3.524'
$ws.Range("E17").Value = $q + '0.57%'
$ws.Range("D18").Value = $q + '2.186'
$ws.Range("E18").Value = $q + '-3.69%'
$ws.Range("E19").Value = $q + '-0.61%'
$ws.Range("D20").Value = $q + '0.03379'
$ws.Range("E20").Value = $q + '3.06%'
$ws.Range("D21").Value = $q + '0.1324'
$ws.Range("E21").Value = $q + '1.23%'
$ws.Range("D22").Value = $q + '3.486'
$ws.Range("E22").Value = $q + '-3.35%'

# Rows 23 and 24 swapped places: ZBToken <-> CoinExToken
$ws.Range("B23").Value = $q + 'CoinExToken'
$ws.Range("C23").Value = $q + 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").Value = $q + '0.04128'
$ws.Range("E23").Value = $q + '-0.14%'
$ws.Range("B24").Value = $q + 'ZBToken'
$ws.Range("C24").Value = $q + 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D24").Value = $q + '0.1410'
$ws.Range("E24").Value = $q + '2.23%'
$ws.Range("D25").Value = $q + '0.001225'
$ws.Range("D26").Value = $q + '0.004149'
$ws.Range("E26").Value = $q + '-7.89%'
$ws.Range("E27").Value = $q + '-0.81%'
$ws.Range("E28").Value = $q + '4.54%'
$ws.Range("E40").Value = $q + '-0.57%'
$ws.Range("D41").Value = $q + '0.1070'
$ws.Range("E41").Value = $q + '-0.22%'
$ws.Range("D42").Value = $q + '0.002200'
$ws.Range("E42").Value = $q + '0.02%'
$ws.Range("E43").Value = $q + '-48.33%'
$ws.Range("D44").Value = $q + '0.008942'
$ws.Range("E44").Value = $q + '-10.95%'
$ws.Range("D45").Value = $q + '0.00005486'
$ws.Range("E45").Value = $q + '7.10%'
$ws.Range("E46").Value = $q + '0.02%'
$ws.Range("E47").Value = $q + '-20.21%'
$ws.Range("D48").Value = $q + '0.002467'
$ws.Range("E48").Value = $q + '-10.83%'
$ws.Range("D49").Value = $q + '0.00002100'
$ws.Range("E49").Value = $q + '0.02%'
$ws.Range("D50").Value = $q + '0.0002000'
$ws.Range("E50").Value = $q + '0.02%'
